# Add team Win/Loss/Tie record columns (AD/AE/AF) to the player data sheet.
# Mirrors the commit "Added team record to data": W/L/T live on the same
# sheet as the player stats instead of a separate sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -----------------------------------------------
# Clone the formatting of the last existing header cell (AC1 - bold,
# centered, thin-bordered) onto the three new header cells, then set
# their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (rows 2-48) ---------------------------------------------
# Every player row gets the same 2009 team record: 78 wins, 84 losses,
# 0 ties.
$ws.Range("AD2:AD48").Value = 78
$ws.Range("AE2:AE48").Value = 84
$ws.Range("AF2:AF48").Value = 0
